# Swap the values of columns E (codeforiati:category-name) and F (codeforiati:group-code)
# for every row in the sheet, including the header row. This matches the upstream
# SectorGroup.xlsx data refresh where the "category-name" and "group-code" columns
# were transposed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    # Force text number format so numeric-looking codes (e.g. "110") are not
    # reinterpreted as numbers when written back.
    $eCell.NumberFormat = "@"
    $fCell.NumberFormat = "@"

    $eCell.Value = $fVal
    $fCell.Value = $eVal
}

Write-Host "Swapped columns E and F for" $lastRow "rows"
